# Apply the "Template_Questions.xlsx" content update:
#  - Rewording of the "Type de réponse attendu" and "Créer autant de colonnes"
#    cells on the "Questions" sheet.
#  - Clear the French button-label values (Valider / Arrêter le questionnaire /
#    Continuer / Démarrer) from the "Textes" sheet (column B).
#  - Clear the long RGPD / presentation descriptive texts from the
#    "Introduction" sheet (column B), keeping the row labels in column A.
#  - Make "Questions" the active/selected sheet instead of "Introduction".

$wb = $excel.ActiveWorkbook

# --- Questions sheet -------------------------------------------------
$wsQuestions = $wb.Worksheets.Item("Questions")

$wsQuestions.Range("B1").Value = 'Type de réponse attendu : "text", "radio" (choix unique), "checkbox" (choix multiple). Si vous souhaitez indiquer un champ "Autre" à votre réponse à choix multiple ou unique, indiquer : "checkbox, text" ou "radio, text"'
$wsQuestions.Range("D1").Value = 'Créer autant de colonnes que de réponses possibles (avec leurs noms). Attention si vous avez un champ "Autre" ne l''ajoutez pas en tant que colonne.'

# --- Textes sheet ------------------------------------------------------
# The button-label cells (B1:B4) are removed outright (content + format),
# not just blanked, so use Clear() rather than ClearContents() here.
$wsTextes = $wb.Worksheets.Item("Textes")

$wsTextes.Range("B1").Clear()
$wsTextes.Range("B2").Clear()
$wsTextes.Range("B3").Clear()
$wsTextes.Range("B4").Clear()

# --- Introduction sheet -------------------------------------------------
$wsIntro = $wb.Worksheets.Item("Introduction")

$wsIntro.Range("B1").ClearContents()
$wsIntro.Range("B2").ClearContents()
$wsIntro.Range("B3").ClearContents()

# --- Selections / active sheet -----------------------------------------
# Each sheet remembers its own last selection; activate before selecting.
$wsTextes.Activate()
$wsTextes.Range("B1:B4").Select()

$wsIntro.Activate()
$wsIntro.Range("B1:B3").Select()

# "Questions" ends up the active/visible tab (was "Introduction" before).
$wsQuestions.Activate()
$wsQuestions.Range("A1:XFD1").Select()
